$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 4 de Octubre de 2020 a las 17:38"

# Each entry: row number, Pais, Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes
$rowUpdates = @(
    @(4, "Estados Unidos", 7606854, 6008, 4828654, 2563870, 0, 53, 214330),
    @(5, "India", 6573678, 26265, 5527934, 943747, 0, 185, 101997),
    @(21, "Italia", 325329, 2578, 231914, 57429, 0, 18, 35986),
    @(22, "Turquia", 323014, 0, 283868, 30762, 0, 0, 8384),
    @(26, "Alemania", 300478, 450, 261900, 28981, 0, 0, 9597),
    @(29, "Canada", 165037, 566, 139401, 16168, 0, 6, 9468),
    @(38, "Republica Dominicana", 114480, 554, 90366, 21980, 0, 6, 2134),
    @(46, "Guatemala", 93963, 215, 82454, 8216, 0, 8, 3293),
    @(49, "Japon", 85339, 571, 78315, 5427, 0, 7, 1597),
    @(61, "Moldavia", 56579, 691, 40955, 14258, 0, 13, 1366),
    @(92, "Jordania", 15640, 891, 5240, 10299, 0, 13, 101),
    @(93, "Senegal", 15094, 26, 12805, 1977, 0, 0, 312),
    @(94, "Zambia", 14974, 0, 14094, 547, 0, 0, 333),
    @(96, "Albania", 14266, 149, 8675, 5195, 0, 4, 396),
    @(100, "Montenegro", 12127, 193, 7953, 3995, 0, 5, 179),
    @(109, "Luxemburgo", 8890, 93, 7428, 1336, 0, 1, 126),
    @(110, "Haiti", 8811, 0, 6949, 1633, 0, 0, 229),
    @(111, "Uganda", 8808, 146, 4736, 3991, 0, 2, 81),
    @(116, "Jamaica", 6895, 100, 2460, 4315, 0, 1, 120),
    @(119, "Cuba", 5809, 29, 5145, 542, 0, 0, 122),
    @(120, "Malaui", 5783, 0, 4539, 1065, 0, 0, 179),
    @(143, "Sri Lanka", 3402, 7, 3258, 131, 0, 0, 13),
    @(162, "Republica de Chipre", 1824, 13, 1369, 433, 0, 0, 22)
)

foreach ($row in $rowUpdates) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
